# Weekly update to the "Hortaliza, Vega Modelo de Temuco - Coliflor" sheet:
# a new daily price record is inserted at row 153 (shifting every existing
# record from row 153 downward by one row), growing the used range from
# A1:R225 to A1:R226.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 153, pushing rows 153:225 down to 154:226.
$ws.Rows.Item(153).Insert(-4121)   # -4121 = xlShiftDown

# Populate the newly inserted row 153 with the new record's data.
$ws.Cells.Item(153, 1).Value  = 10
$ws.Cells.Item(153, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(153, 3).Value  = "La Araucanía"
$ws.Cells.Item(153, 4).Value  = 44460
$ws.Cells.Item(153, 5).Value  = 9
$ws.Cells.Item(153, 6).Value  = 100112008
$ws.Cells.Item(153, 7).Value  = "Coliflor"
$ws.Cells.Item(153, 8).Value  = "Sin especificar"
$ws.Cells.Item(153, 9).Value  = "Primera"
$ws.Cells.Item(153, 10).Value = 1000
$ws.Cells.Item(153, 11).Value = 1100
$ws.Cells.Item(153, 12).Value = 1100
$ws.Cells.Item(153, 13).Value = 1100
$ws.Cells.Item(153, 14).Value = "$/unidad"
$ws.Cells.Item(153, 15).Value = "Región Metropolitana"
$ws.Cells.Item(153, 16).Value = 1100
$ws.Cells.Item(153, 17).Value = 1
$ws.Cells.Item(153, 18).Value = "Hortaliza"
